# Senior Project Tracking RTC.xlsx
# Add a new log entry (row 8) recording work on the "dummy system" that
# pulls data from RTC into the Task Management System.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 8: Date / Start time / End time / Hours (formula) / Reason ---
$ws.Range("A8").Value = 42412
$ws.Range("A8").NumberFormat = "mm-dd-yy"

$ws.Range("B8").Value = 0.4375
$ws.Range("B8").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Range("C8").Value = 0.77083333333333337
$ws.Range("C8").NumberFormat = $ws.Range("C2").NumberFormat

$ws.Range("D8").Formula = "=MOD(IF(ISBLANK(C8),B8, C8)-B8, 1)*24"
$ws.Range("D8").NumberFormat = $ws.Range("D2").NumberFormat

$ws.Range("E8").Value = "Began experimenting with dummy system to pull data from RTC into the Task Management System."

# Column A now holds a real date value (42412) alongside the existing
# "d-mmm" formatted dates, so widen/autofit it to show the new format.
$ws.Columns.Item(1).AutoFit()

# Move the active selection down to the next empty row, as Excel does
# after the last data entry.
$ws.Range("E9").Select()
